$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column J
$ws.Range("J1").Value = "decimal[decimal]"

# Data rows
$ws.Range("J2").Value = 15.3
$ws.Range("J3").Value = 14.3

# J3 gets an explicit "General" number format applied (matches xf with applyNumberFormat)
$ws.Range("J3").NumberFormat = "General"

# Update selection to match post-edit state
$ws.Range("J3").Select()
